$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.400.74"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "2.062.89"
$ws.Range("E3").Value = "  +4.97%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.25"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +9.26%  "
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0761"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  +4.21%  "
$ws.Range("D13").Value = "2.370.45"
$ws.Range("E13").Value = "  +5.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.35"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.15"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.778"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").Value = "2.083.69"
$ws.Range("E18").Value = "  +6.13%  "
$ws.Range("D19").Value = "37.460.42"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +24.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.94"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.05%  "
$ws.Range("D22").Value = "0.0₃0812"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "224.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.22%  "
$ws.Range("E26").Value = "  +2.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.126"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0627"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +14.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.49%  "
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +14.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.62"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +27.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0980"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +12.31%  "
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.484.34"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.85"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +11.87%  "
$ws.Range("E46").Value = "  +6.55%  "
$ws.Range("E47").Value = "  +7.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.88"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.95"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.47%  "
